# ------------------------------------------------------------------
# Add the new "upvote_novote" sheet (goes after baseline_novote and
# becomes the active tab), and populate it with the upvote-vs-novote
# correlation table.
# ------------------------------------------------------------------
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "upvote_novote"

# ---- Header row (row 1) ----
# Seed formats by copying the matching header cells from baseline_upvote
# so the number formats / bold styling line up with the rest of the
# workbook instead of minting near-duplicate custom formats.
$ws1.Range("A1").Copy($ws.Range("A1"))
$ws1.Range("B1").Copy($ws.Range("B1"))
$ws1.Range("D1").Copy($ws.Range("C1"))
$ws1.Range("E1").Copy($ws.Range("E1"))
$ws1.Range("F1").Copy($ws.Range("F1"))
$ws1.Range("G1").Copy($ws.Range("G1"))
$ws1.Range("A1").Copy($ws.Range("D1"))
$ws1.Range("A1").Copy($ws.Range("G1"))

$ws.Range("A1").Value = "treatment subreddit"
$ws.Range("B1").Value = "control subreddit"
$ws.Range("C1").Value = "upvote only dates"
$ws.Range("D1").Value = "no vote dates"
$ws.Range("E1").Value = "upvote only n"
$ws.Range("F1").Value = "novote n"
$ws.Range("G1").Value = "correlation"

# "novote n" doesn't already exist as a rich string anywhere in the
# workbook, so split the trailing "n" into its own bold+italic run
# (matches "baseline n" / "upvote only n" / "no votes n" elsewhere).
$fChars = $ws.Range("F1").Characters(8, 1)
$fChars.Font.Bold = $true
$fChars.Font.Italic = $true

# ---- \Conservative block (rows 2-5) ----
$ws1.Range("E2").Copy($ws.Range("E2"))
$ws1.Range("F2").Copy($ws.Range("F2"))
$ws.Range("A2").Value = "\Conservative"
$ws.Range("C2").Value = "2013-11-01 to 2014-06-30"
$ws.Range("D2").Value = "2014-07-01 to 2018-01-01"
$ws.Range("E2").Value = 131941
$ws.Range("F2").Value = 1535763

$ws.Range("B3").Value = "\progressive"
$ws1.Range("E3").Copy($ws.Range("E3"))
$ws1.Range("F3").Copy($ws.Range("F3"))
$ws.Range("E3").Value = 24238
$ws.Range("F3").Value = 131108
$ws.Range("G3").Value = 0.998

$ws.Range("B4").Value = "\Liberal"
$ws.Range("E4").Value = 19745
$ws.Range("F4").Value = 96124
$ws.Range("E4").NumberFormat = "#,##0"
$ws.Range("F4").NumberFormat = "#,##0"
$ws.Range("G4").Value = 0.997

$ws.Range("B5").Value = "\Republican"
$ws.Range("E5").Value = 14449
$ws.Range("F5").Value = 205248
$ws.Range("E5").NumberFormat = "#,##0"
$ws.Range("F5").NumberFormat = "#,##0"
$ws.Range("G5").Value = 0.992

# ---- \GenderCritical block (rows 7-8) ----
$ws1.Range("C18").Copy($ws.Range("C7"))
$ws.Range("A7").Value = "\GenderCritical"
$ws.Range("C7").Value = "2013-10-03 to 2014-06-30"
$ws.Range("D7").Value = "2014-07-01 to 2018-04-01"
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = 7857
$ws.Range("F7").Value = 429107
$ws.Range("E7").NumberFormat = "#,##0"
$ws.Range("F7").NumberFormat = "#,##0"
$ws.Range("7:7").RowHeight = 17

$ws.Range("B8").Value = "\iran"
$ws.Range("E8").Value = 15035
$ws.Range("F8").Value = 123372
$ws.Range("E8").NumberFormat = "#,##0"
$ws.Range("F8").NumberFormat = "#,##0"
$ws.Range("G8").Value = 0.918

# ---- \politics block (rows 10-16) ----
$ws.Range("A10").Value = "\politics"
$ws.Range("C10").Value = "2014-01-30 to 2014-12-29"
$ws.Range("D10").Value = "2014-12-30 to 2018-04-01"
$ws.Range("E10").Value = 2070909
$ws.Range("F10").Value = 50203736
$ws.Range("E10").NumberFormat = "#,##0"
$ws.Range("F10").NumberFormat = "#,##0"

$ws.Range("B11").Value = "\Anarcho_Capitalism"
$ws.Range("E11").Value = 278502
$ws.Range("F11").Value = 931707
$ws.Range("E11").NumberFormat = "#,##0"
$ws.Range("F11").NumberFormat = "#,##0"
$ws.Range("G11").Value = 0.996

$ws.Range("B12").Value = "\privacy"
$ws.Range("E12").Value = 40575
$ws.Range("F12").Value = 295572
$ws.Range("E12").NumberFormat = "#,##0"
$ws.Range("F12").NumberFormat = "#,##0"
$ws.Range("G12").Value = 0.99

$ws.Range("B13").Value = "\conspiracy"
$ws.Range("E13").Value = 756534
$ws.Range("F13").Value = 3420420
$ws.Range("E13").NumberFormat = "#,##0"
$ws.Range("F13").NumberFormat = "#,##0"
$ws.Range("G13").Value = 0.985

$ws.Range("B14").Value = "\PoliticalDiscussion"
$ws.Range("E14").Value = 253573
$ws.Range("F14").Value = 791679
$ws.Range("E14").NumberFormat = "#,##0"
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("G14").Value = 0.977

$ws.Range("B15").Value = "\MensRights"
$ws.Range("E15").Value = 477218
$ws.Range("F15").Value = 1474394
$ws.Range("E15").NumberFormat = "#,##0"
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 0.973

$ws.Range("B16").Value = "\PoliticalHumor"
$ws.Range("E16").Value = 27013
$ws.Range("F16").Value = 1645018
$ws.Range("E16").NumberFormat = "#,##0"
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 0.968

# ---- Column widths (match the "best fit" sizing used on other sheets) ----
$ws.Columns("A").ColumnWidth = 17.83203125
$ws.Columns("B").ColumnWidth = 18.5
$ws.Columns("C").ColumnWidth = 22.83203125
$ws.Columns("D").ColumnWidth = 22.83203125
$ws.Columns("E").ColumnWidth = 15
$ws.Columns("F").ColumnWidth = 13.6640625

$ws.Range("G17").Select()
